$wb = $excel.ActiveWorkbook
$results = $wb.Worksheets.Item(1)
$players = $wb.Worksheets.Item(2)

# ============================================================
# Sheet "results": convert the time column to free-text times
# and add the missing third race result (Mike vs Aidan B).
# ============================================================

# Existing row 2 (Justine / Jake / Donkey Kong) - time becomes text
$results.Range("D1").NumberFormat = "@"
$results.Range("D2").NumberFormat = "@"
$results.Range("D2").Value2 = "03:12.670"

# New row 3 - Mike / Aidan B / Toad
$results.Range("A3").Value2 = "Mike"
$results.Range("B3").Value2 = "Aidan B"
$results.Range("C3").Value2 = "Toad"
$results.Range("D3").NumberFormat = "@"
$results.Range("D3").Value2 = "03:00.809"
$results.Range("E3").Value2 = 45978

# Column widths for D (time) and E (date), matching the new wider/text layout
$results.Columns("D").ColumnWidth = 11.4
$results.Columns("E").ColumnWidth = 19.1

$results.Activate()
$results.Range("D3").Select()

# ============================================================
# Sheet "players": add the two new players, Mike and Aidan B
# ============================================================

$players.Range("C3:D3").Copy()
$players.Range("C4:D4").PasteSpecial(-4122)
$players.Range("C5:D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$players.Range("A4").Value2 = "Mike"
$players.Range("B4").Value2 = "mike.jpg"
$players.Range("C4").Value2 = $players.Range("C3").Value2
$players.Range("D4").Value2 = $players.Range("D3").Value2

$players.Range("A5").Value2 = "Aidan B"
$players.Range("B5").Value2 = "aidanb.jpg"
$players.Range("C5").Value2 = $players.Range("C3").Value2
$players.Range("D5").Value2 = $players.Range("D3").Value2

$players.Rows(4).RowHeight = $players.Rows(3).RowHeight
$players.Rows(5).RowHeight = $players.Rows(3).RowHeight

$players.Columns("B").ColumnWidth = 9.8

$players.Activate()
$players.Range("B4").Select()

$results.Activate()
